# Update odds in the "Jogos do Dia Betfair Back Lay" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.02
$ws.Range("H2").Value = 1.02
$ws.Range("N2").Value = 1.3
$ws.Range("P2").Value = 1.3
$ws.Range("T2").Value = 1.03
$ws.Range("U2").Value = 1.03

# Row 3
$ws.Range("F3").Value = 1.37
$ws.Range("G3").Value = 1.39
$ws.Range("H3").Value = 10.5
$ws.Range("I3").Value = 14
$ws.Range("J3").Value = 5.1
$ws.Range("K3").Value = 5.6
$ws.Range("N3").Value = 4
$ws.Range("P3").Value = 2.02
$ws.Range("Q3").Value = 1.84

# Row 4
$ws.Range("G4").Value = 4.4
$ws.Range("H4").Value = 2.1
$ws.Range("K4").Value = 3.5
$ws.Range("P4").Value = 1.67
$ws.Range("Q4").Value = 2.28

# Row 5
$ws.Range("F5").Value = 2.16
$ws.Range("G5").Value = 2.28
$ws.Range("H5").Value = 3.55
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 3.35
$ws.Range("K5").Value = 3.6
$ws.Range("P5").Value = 1.76
$ws.Range("Q5").Value = 2.12

# Row 6
$ws.Range("F6").Value = 2.48
$ws.Range("G6").Value = 2.56
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 3.65
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 3.35
$ws.Range("P6").Value = 1.58
$ws.Range("Q6").Value = 2.44

# Row 7
$ws.Range("F7").Value = 2.46
$ws.Range("G7").Value = 2.7
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 3.85
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 3.35
$ws.Range("P7").Value = 1.6
$ws.Range("Q7").Value = 2.42

$wb.Save()
